$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.215.04"
$ws.Range("E2").Value = "'  +0.74%  "

# Row 3
$ws.Range("D3").Value = "'1.852.21"
$ws.Range("E3").Value = "'  +1.19%  "

# Row 4
$ws.Range("E4").Value = "'  -0.32%  "

# Row 5
$ws.Range("D5").Value = "'313.60"
$ws.Range("E5").Value = "'  +0.37%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  -0.25%  "

# Row 7
$ws.Range("D7").Value = "'0.4606"
$ws.Range("E7").Value = "'  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.3705"
$ws.Range("E8").Value = "'  -0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.07284"
$ws.Range("E9").Value = "'  -0.89%  "

# Row 10
$ws.Range("D10").Value = "'0.8841"
$ws.Range("E10").Value = "'  +1.06%  "

# Row 11
$ws.Range("D11").Value = "'20.02"
$ws.Range("E11").Value = "'  +1.07%  "

# Row 12
$ws.Range("D12").Value = "'0.07829"
$ws.Range("E12").Value = "'  -1.49%  "

# Row 13
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.906.19"
$ws.Range("E13").Value = "'  +5.76%  "

# Row 14
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.374"
$ws.Range("E14").Value = "'  +0.56%  "

# Row 15
$ws.Range("D15").Value = "'6.511"
$ws.Range("E15").Value = "'  -0.53%  "

# Row 16
$ws.Range("D16").Value = "'91.38"
$ws.Range("E16").Value = "'  -0.15%  "

# Row 17
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "'  -0.32%  "

# Row 18
$ws.Range("D18").Value = "'0.000008922"
$ws.Range("E18").Value = "'  +0.14%  "

# Row 19
$ws.Range("E19").Value = "'  -0.04%  "

# Row 20
$ws.Range("D20").Value = "'14.71"
$ws.Range("E20").Value = "'  -0.51%  "

# Row 21
$ws.Range("D21").Value = "'27.245.87"
$ws.Range("E21").Value = "'  +1.55%  "

# Row 22
$ws.Range("D22").Value = "'5.101"
$ws.Range("E22").Value = "'  -0.29%  "

# Row 23
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "'  -0.55%  "

# Row 24
$ws.Range("D24").Value = "'2.087.03"
$ws.Range("E24").Value = "'  +4.61%  "

# Row 25
$ws.Range("D25").Value = "'1.930"
$ws.Range("E25").Value = "'  +4.98%  "

# Row 26
$ws.Range("D26").Value = "'151.54"
$ws.Range("E26").Value = "'  -1.02%  "

# Row 27
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "'  -0.26%  "

# Row 28
$ws.Range("D28").Value = "'2.062"
$ws.Range("E28").Value = "'  +0.64%  "

# Row 29
$ws.Range("D29").Value = "'115.77"
$ws.Range("E29").Value = "'  +0.11%  "

# Row 30
$ws.Range("D30").Value = "'5.046"
$ws.Range("E30").Value = "'  -2.22%  "

# Row 31
$ws.Range("D31").Value = "'0.08821"
$ws.Range("E31").Value = "'  -0.96%  "

# Row 32
$ws.Range("D32").Value = "'3.095"
$ws.Range("E32").Value = "'  +4.44%  "

# Row 33
$ws.Range("D33").Value = "'0.7624"
$ws.Range("E33").Value = "'  +4.13%  "

# Row 34
$ws.Range("E34").Value = "'  +3.47%  "

# Row 35
$ws.Range("E35").Value = "'  +1.54%  "

# Row 36
$ws.Range("D36").Value = "'2.734"
$ws.Range("E36").Value = "'  +10.25%  "

# Row 37
$ws.Range("D37").Value = "'1.084"
$ws.Range("E37").Value = "'  +1.18%  "

# Row 38
$ws.Range("D38").Value = "'0.01946"
$ws.Range("E38").Value = "'  -0.16%  "

# Row 39
$ws.Range("D39").Value = "'0.05240"
$ws.Range("E39").Value = "'  +0.01%  "

# Row 40
$ws.Range("D40").Value = "'2.949"
$ws.Range("E40").Value = "'  +0.35%  "

# Row 41
$ws.Range("D41").Value = "'7.054"
$ws.Range("E41").Value = "'  -0.91%  "

# Row 42
$ws.Range("D42").Value = "'0.5100"
$ws.Range("E42").Value = "'  -1.08%  "

# Row 43
$ws.Range("D43").Value = "'0.1624"
$ws.Range("E43").Value = "'  -0.21%  "

# Row 44
$ws.Range("D44").Value = "'8.363"
$ws.Range("E44").Value = "'  +1.82%  "

# Row 45
$ws.Range("D45").Value = "'0.4787"
$ws.Range("E45").Value = "'  -1.06%  "

# Row 46
$ws.Range("D46").Value = "'10.31"
$ws.Range("E46").Value = "'  +1.29%  "

# Row 47
$ws.Range("E47").Value = "'  -0.38%  "

# Row 48
$ws.Range("D48").Value = "'102.07"
$ws.Range("E48").Value = "'  -0.20%  "

# Row 49
$ws.Range("D49").Value = "'1.638"
$ws.Range("E49").Value = "'  +0.19%  "

# Row 50
$ws.Range("D50").Value = "'0.06207"
$ws.Range("E50").Value = "'  +0.12%  "

# Row 51
$ws.Range("D51").Value = "'65.44"
$ws.Range("E51").Value = "'  +0.83%  "
